$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at position 8 (this pushes the old row 8 -- the
#    totals row -- down to row 9, and the old row 9 -- the footer row
#    with date / page / developed-by -- down to row 10). Excel auto
#    shifts the merged cell ranges that live below the insertion point.
# ------------------------------------------------------------------
$ws.Rows("8:8").Insert()

# ------------------------------------------------------------------
# 2. Populate the new item row (row 8) with the new medicine sold,
#    copying the cell formatting from row 7 (the first item row) so
#    the new row looks identical in style to the existing item rows.
# ------------------------------------------------------------------
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "LASIX 40MG 24 TAB."
$ws.Range("H8").Value = "2:0"

# L8 holds the quantity as text (matches L7's shared-string "1"), so
# force text formatting before assigning, then restore the original
# number-format/style that row 7 uses for that column.
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "1"
$ws.Range("L7").Copy()
$ws.Range("L8").PasteSpecial(-4122)

$ws.Range("N8").Value = "30.00"
$ws.Range("P8").Value = "15.0000"
$ws.Range("Q8").Value = "0:1"

# Re-apply row 7's styling one more time over the whole row so every
# cell (including the ones just written) keeps the same per-column
# look (border / fill / font / number format) as the template row.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)

# Put the real values back (PasteSpecial above only touches formats,
# but do this defensively in case a prior PasteSpecial altered them).
$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "LASIX 40MG 24 TAB."
$ws.Range("H8").Value = "2:0"
$ws.Range("N8").Value = "30.00"
$ws.Range("P8").Value = "15.0000"
$ws.Range("Q8").Value = "0:1"

$ws.Application.CutCopyMode = 0

# Merge the cell groups for the new row, mirroring row 7's layout.
$ws.Range("A8:B8").Merge()
$ws.Range("C8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()
$ws.Range("N8:O8").Merge()

# ------------------------------------------------------------------
# 3. Row heights: the new row 8 takes on the height the totals row
#    used to have (24.75), while the totals row -- now row 9 -- takes
#    on the height the first item row (row 7) uses (25.5). Row 10
#    (old row 9, the footer) keeps its original 16.5 height.
# ------------------------------------------------------------------
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5

# ------------------------------------------------------------------
# 4. Update the totals row (now row 9): the grand total in column P
#    must include the price of the newly added item.
# ------------------------------------------------------------------
$ws.Range("P9").Value = 16.829999999999998 + 15.0000
$ws.Range("Q9").Value = ""

# ------------------------------------------------------------------
# 5. Update the footer timestamp (now row 10) to the new save time.
# ------------------------------------------------------------------
$ws.Range("A10").Value = "Friday, 18 July, 2025 3:59 PM"

Write-Host "edit complete"
